# March Deskcount.xlsx update
# - Flip "Include in Occupancy Calculation" (column F) from Yes to No for a
#   handful of offices (Los Angeles, Orlando, Philadelphia, Tampa, Dublin,
#   Santiago, Sao Paulo).
# - Correct the Melbourne deskcount (C45) from 30 to 32.
# - Move the active selection to D44 (mirrors the author's last selection
#   when saving, with the view scrolled down near row 18).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deskcount")

# Rows whose "Include in Occupancy Calculation" flag moved from Yes to No.
$noRows = @(20, 27, 28, 39, 41, 48, 49)
foreach ($r in $noRows) {
    $ws.Cells.Item($r, 6).Value = "No"
}

# Melbourne's deskcount correction.
$ws.Range("C45").Value = 32

# Restore the author's on-save selection / scroll position.
$ws.Activate()
$ws.Range("A18").Select()
$ws.Range("D44").Select()
